# Update countries & provincias Spain
# Applies the 21-Abril-2020 19:22 data refresh to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp
#  - updates case counters for several countries (rows 4, 8, 30, 87)
#  - Jordania overtakes Taiwan in the ranking (rows 107/108 swap + new Jordania numbers)
#  - Guinea Ecuatorial overtakes Maldivas & Liechtenstein (rows 148/149/150 rotate + new Guinea Ecuatorial numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 19:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 803349
$ws.Range("C4").Value = 10590
$ws.Range("E4").Value = 684397
$ws.Range("G4").Value = 1121
$ws.Range("H4").Value = 43635

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 148007
$ws.Range("C8").Value = 942
$ws.Range("E8").Value = 47868
$ws.Range("G8").Value = 77
$ws.Range("H8").Value = 4939

# --- Ecuador (row 30) ---
$ws.Range("B30").Value = 10398
$ws.Range("C30").Value = 270
$ws.Range("D30").Value = 1207
$ws.Range("E30").Value = 8671
$ws.Range("F30").Value = 137
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 520

# --- Tunez (row 87) ---
$ws.Range("F87").Value = 35

# --- Jordania / Taiwan swap ranking positions (rows 107-108) ---
# Row 107 becomes Jordania with refreshed totals, row 108 becomes Taiwan (unchanged totals)
$ws.Range("A107").Value = "Jordania"
$ws.Range("B107").Value = 428
$ws.Range("C107").Value = 3
$ws.Range("D107").Value = 297
$ws.Range("E107").Value = 124
$ws.Range("F107").Value = 5
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 7

$ws.Range("A108").Value = "Taiwan"
$ws.Range("B108").Value = 425
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 217
$ws.Range("E108").Value = 202
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 6

# --- Guinea Ecuatorial overtakes Maldivas & Liechtenstein (rows 148-150) ---
# Row 148 becomes Guinea Ecuatorial with refreshed totals,
# row 149 becomes Maldivas (unchanged totals), row 150 becomes Liechtenstein (unchanged totals)
$ws.Range("A148").Value = "Guinea Ecuatorial"
$ws.Range("B148").Value = 83
$ws.Range("C148").Value = 4
$ws.Range("D148").Value = 7
$ws.Range("E148").Value = 76
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0

$ws.Range("A149").Value = "Maldivas"
$ws.Range("B149").Value = 83
$ws.Range("C149").Value = 14
$ws.Range("D149").Value = 16
$ws.Range("E149").Value = 67
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 0

$ws.Range("A150").Value = "Liechtenstein"
$ws.Range("B150").Value = 81
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 55
$ws.Range("E150").Value = 25
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 1

$wb.Save()
